{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" site-chrome\n// paragraphs (plus the blank paragraph that separated them from the\n// preceding \"Requisitos\" text), mirroring the Jekyll site rebuild that\n// dropped this boilerplate footer from the generated page.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two text paragraphs that must go.\nlet jupiterIdx = -1;\nlet copyrightIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (jupiterIdx === -1 && t.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIdx = i;\n  }\n  if (copyrightIdx === -1 && t.indexOf(\"Powered by Jekyll\") !== -1) {\n    copyrightIdx = i;\n  }\n}\n\nif (jupiterIdx === -1 || copyrightIdx === -1) {\n  throw new Error(\"Could not locate the footer paragraphs to remove.\");\n}\n\n// The blank separator paragraph directly above the \"Ver no Jupiter\" line\n// (inserted right after the \"Requisitos\" answer paragraph) is removed too.\nconst blankIdx = jupiterIdx - 1;\n\n// Delete from the bottom up so earlier indices stay valid.\nitems[copyrightIdx].delete();\nitems[jupiterIdx].delete();\nif (blankIdx >= 0 && items[blankIdx].text === \"\") {\n  items[blankIdx].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" site-chrome\n# paragraphs (plus the blank paragraph that separated them from the\n# preceding \"Requisitos\" answer), mirroring the Jekyll site rebuild that\n# dropped this boilerplate footer from the generated page.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$jupiterIdx = -1\n$copyrightIdx = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($jupiterIdx -eq -1 -and $t -like \"*Ver no Jupiter*\") {\n        $jupiterIdx = $i\n    }\n    if ($copyrightIdx -eq -1 -and $t -like \"*Powered by Jekyll*\") {\n        $copyrightIdx = $i\n    }\n}\n\nif ($jupiterIdx -eq -1 -or $copyrightIdx -eq -1) {\n    throw \"Could not locate the footer paragraphs to remove.\"\n}\n\n# The blank separator paragraph directly above the \"Ver no Jupiter\" line.\n$blankIdx = $jupiterIdx - 1\n\n# Delete highest index first so the lower indices stay valid.\n$d.Paragraphs.Item($copyrightIdx).Range.Delete()\n$d.Paragraphs.Item($jupiterIdx).Range.Delete()\nif ($blankIdx -ge 1 -and $d.Paragraphs.Item($blankIdx).Range.Text.Trim() -eq \"\") {\n    $d.Paragraphs.Item($blankIdx).Range.Delete()\n}\n"}
